# Update "想去人数" (want-to-go count) values in column F for rows 3, 6, 7, 8, 9
# on both the "展览" and "全部类型" worksheets, matching newly scraped totals.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3 = 497
    6 = 46
    7 = 40
    8 = 2009
    9 = 4088
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
